$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some new "Price" values look like plain numbers (e.g. "236.23", "1.000").
# Typing those with the default General format would make Excel store them
# as floating point numbers (losing trailing zeros / exact text layout), so
# those specific cells are switched to a Text number format first -- exactly
# like a user would do in the UI -- before the literal value is entered.

$ws.Range("D2").Value = "29.182.53"
$ws.Range("E2").Value = "  -0.54%  "
$ws.Range("D3").Value = "1.827.54"
$ws.Range("E3").Value = "  -0.72%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9991"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.23"
$ws.Range("E5").Value = "  -1.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6044"
$ws.Range("E6").Value = "  -3.91%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07119"
$ws.Range("E8").Value = "  -4.58%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2818"
$ws.Range("E9").Value = "  -2.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.03"
$ws.Range("E10").Value = "  -3.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07653"
$ws.Range("D12").Value = "1.840.76"
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.772"
$ws.Range("E13").Value = "  -4.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6408"
$ws.Range("E14").Value = "  -5.47%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000009948"
$ws.Range("E15").Value = "  -2.40%  "
$ws.Range("D16").Value = "2.073.20"
$ws.Range("E16").Value = "  -0.78%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "79.47"
$ws.Range("E17").Value = "  -3.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.996"
$ws.Range("E18").Value = "  -3.87%  "
$ws.Range("D19").Value = "29.206.76"
$ws.Range("E19").Value = "  -0.43%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "231.72"
$ws.Range("E20").Value = "  +1.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.78"
$ws.Range("E22").Value = "  -4.40%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.049"
$ws.Range("E23").Value = "  -4.94%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.002"
$ws.Range("E24").Value = "  +0.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "155.24"
$ws.Range("E25").Value = "  -2.23%  "
$ws.Range("E26").Value = "  -4.95%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1283"
$ws.Range("E27").Value = "  -5.26%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.67"
$ws.Range("E28").Value = "  -4.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06875"
$ws.Range("E29").Value = "  +6.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.458"
$ws.Range("E30").Value = "  +0.82%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.451"
$ws.Range("E31").Value = "  -2.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.840"
$ws.Range("E32").Value = "  -5.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.804"
$ws.Range("E33").Value = "  -6.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.136"
$ws.Range("E34").Value = "  -0.17%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.721"
$ws.Range("E35").Value = "  -6.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6629"
$ws.Range("E36").Value = "  -4.30%  "
$ws.Range("E37").Value = "  -1.56%  "
$ws.Range("D38").Value = "1.234.21"
$ws.Range("E38").Value = "  -0.50%  "
$ws.Range("E39").Value = "  -2.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01762"
$ws.Range("E40").Value = "  -5.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.598"
$ws.Range("E41").Value = "  -2.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9333"
$ws.Range("E42").Value = "  +0.35%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.001"
$ws.Range("E43").Value = "  +0.20%  "
$ws.Range("D44").Value = "1.995.31"
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "100.10"
$ws.Range("E45").Value = "  -0.70%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "63.47"
$ws.Range("E46").Value = "  -3.25%  "
$ws.Range("E47").Value = "  +1.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.639"
$ws.Range("E48").Value = "  -4.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.580"
$ws.Range("E49").Value = "  -6.64%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05585"
$ws.Range("E50").Value = "  -1.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.487"
$ws.Range("E51").Value = "  -5.96%  "
